# Fruta / hortaliza, semanal
# Inserts a new weekly price record at row 230 (Membrillo / Champion / Primera,
# date serial 45093 = 2023-06-16), which pushes every following record down by
# one row (old row 230 becomes 231, ..., old row 293 becomes 294).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row above the current row 230, shifting rows 230:293
# down to 231:294. Excel preserves the row/column formatting (e.g. the date
# number format on column D) automatically during the insert.
$ws.Rows("230:230").Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Cells.Item(230, 1).Value = 10
$ws.Cells.Item(230, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(230, 3).Value = "La Araucanía"
$ws.Cells.Item(230, 4).Value = 45093
$ws.Cells.Item(230, 5).Value = 9
$ws.Cells.Item(230, 6).Value = "Fruta"
$ws.Cells.Item(230, 7).Value = 100104
$ws.Cells.Item(230, 8).Value = "Frutos de pepita"
$ws.Cells.Item(230, 9).Value = 100104003
$ws.Cells.Item(230, 10).Value = "Membrillo"
$ws.Cells.Item(230, 11).Value = "Champion"
$ws.Cells.Item(230, 12).Value = "Primera"
$ws.Cells.Item(230, 13).Value = 110
$ws.Cells.Item(230, 14).Value = 13000
$ws.Cells.Item(230, 15).Value = 14000
$ws.Cells.Item(230, 16).Value = 13455
$ws.Cells.Item(230, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(230, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(230, 19).Value = 748
$ws.Cells.Item(230, 20).Value = 18
